{"js": "// The \"Referencias (References)\" title paragraph (first table, header row)\n// loses its \" (References)\" suffix (the space, italic \"(\", \"References\",\n// \")\" runs), leaving just \"Referencias\". Word's auto-managed \"_GoBack\"\n// bookmark (which marks the last edited spot) moves from the empty\n// paragraph where it used to sit down near the end of the body to the end\n// of this now-shorter paragraph.\n\nconst body = context.document.body;\n\n// Locate the paragraph that holds \"Referencias (References)\" via search,\n// then grab its enclosing paragraph.\nconst hits = body.search(\"Referencias\", { matchCase: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find 'Referencias' paragraph\");\n}\n\nconst titlePara = hits.items[0].paragraphs.getFirst();\nconst titleRange = titlePara.getRange();\n\n// Remove the \" (References)\" suffix that follows \"Referencias\".\nconst suffixHits = titleRange.search(\" (References)\", { matchCase: false });\nsuffixHits.load(\"items\");\nawait context.sync();\n\nif (suffixHits.items.length > 0) {\n  suffixHits.items[0].delete();\n  await context.sync();\n}\n\n// Move the \"_GoBack\" bookmark: drop it from wherever it currently is, and\n// re-insert it (collapsed) right at the end of the \"Referencias\" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\ntitlePara.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The \"Referencias (References)\" title paragraph (first table, header row)\n# loses its \" (References)\" suffix (the space, italic \"(\", \"References\",\n# \")\" runs), leaving just \"Referencias\". Word's auto-managed \"_GoBack\"\n# bookmark (which marks the last edited spot) moves from the empty\n# paragraph where it used to sit down near the end of the body to the end\n# of this now-shorter paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Referencias (References)\" title and collapse the found range\n# to right after \"Referencias\". Plant the \"_GoBack\" bookmark there BEFORE\n# deleting any text, so it ends up exactly where the real edit would leave\n# the cursor.\n$titleEnd = $d.Content\n$titleEnd.Find.Execute(\"Referencias\")\n$titleEnd.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $titleEnd)\n\n# Delete the \" (References)\" suffix (space + italic \"(References)\") that\n# still follows \"Referencias\".\n$suffix = $d.Content\n$suffix.Find.Execute(\"Referencias\")\n$suffix.Collapse(0)\n$suffix.MoveEnd(1, 13)\n$suffix.Delete()\n"}
